# Add an "Age" column of data to Sheet1 and populate Sheet2 with
# Keyword/TestNg rows, then leave Sheet2 as the active sheet/tab with
# row 3 selected (mirrors the author's "Add files via upload" re-save).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: add numeric "Age" values first so the old shared string
#     ("Pass") in C2 is orphaned and reclaimed by the save step, then
#     write the new "Age" header into the freed slot. ---
$ws1.Range("C2").Value = 36
$ws1.Range("C3").Value = 37
$ws1.Range("C4").Value = 38
$ws1.Range("C1").Value = "Age"

# Move Sheet1's selection off the old A20 cell.
$ws1.Range("A2").Select()

# --- Sheet2: populate with new shared strings (TestNg created before
#     Keyword so the shared-string table order matches the target). ---
$ws2.Range("A2").Value = "TestNg"
$ws2.Range("A1").Value = "Keyword"

# Sheet2 becomes the active/selected sheet, with row 3 selected.
$ws2.Activate()
$ws2.Rows.Item(3).Select()
